$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 456.41
$ws.Range("I15").Value = 456.41
$ws.Range("K15").Value = 1369.23
$ws.Range("M15").Value = -1200.23
$ws.Range("H19").Value = 620
$ws.Range("I19").Value = 608
$ws.Range("J19").Value = 650
$ws.Range("K19").Value = 608
$ws.Range("L19").Value = 650
$ws.Range("M19").Value = -433
$ws.Range("N19").Value = -1000
$ws.Range("H28").Value = 709.64703
$ws.Range("I28").Value = 744.625
$ws.Range("J28").Value = 150
$ws.Range("K28").Value = 744.625
$ws.Range("L28").Value = 150
$ws.Range("M28").Value = -259.625
$ws.Range("N28").Value = -1120
$ws.Range("H127").Value = 1337.2307
$ws.Range("I127").Value = 793.8570999999999
$ws.Range("J127").Value = 1390.8029
$ws.Range("K127").Value = 2381.5713
$ws.Range("L127").Value = 4172.4087
$ws.Range("M127").Value = 2578.4287
$ws.Range("N127").Value = -14092.4087
$ws.Range("H129").Value = 977.9836
$ws.Range("I129").Value = 504.5
$ws.Range("J129").Value = 1070.8235
$ws.Range("K129").Value = 1513.5
$ws.Range("L129").Value = 3212.4705
$ws.Range("M129").Value = 3486.5
$ws.Range("N129").Value = -13212.4705

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = ""
$ws.Range("H61").Value = 2362.75
$ws.Range("I61").Value = 1965.2222
$ws.Range("K61").Value = 1965.2222
$ws.Range("M61").Value = -1753.2222
$ws.Range("H132").Value = 646760.1
$ws.Range("I132").Value = 953140.5600000001
$ws.Range("J132").Value = 3361.2
$ws.Range("K132").Value = 2859421.68
$ws.Range("L132").Value = 10083.6
$ws.Range("M132").Value = -2856891.68
$ws.Range("N132").Value = -15143.6
$ws.Range("H133").Value = 41689
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 41689
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 41689
$ws.Range("M133").Value = ""
$ws.Range("N133").Value = -46749
$ws.Range("H136").Value = 2362.75
$ws.Range("I136").Value = 1965.2222
$ws.Range("K136").Value = 5895.6666
$ws.Range("M136").Value = -3345.6666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2848.8923
$ws.Range("I31").Value = 2219.7715
$ws.Range("J31").Value = 3582.8667
$ws.Range("K31").Value = 2219.7715
$ws.Range("L31").Value = 3582.8667
$ws.Range("M31").Value = -1924.7715
$ws.Range("N31").Value = -4172.8667
$ws.Range("H34").Value = 2848.8923
$ws.Range("I34").Value = 2219.7715
$ws.Range("J34").Value = 3582.8667
$ws.Range("K34").Value = 2219.7715
$ws.Range("L34").Value = 3582.8667
$ws.Range("M34").Value = -2017.7715
$ws.Range("N34").Value = -3986.8667
$ws.Range("H58").Value = 1295.0286
$ws.Range("I58").Value = 1236.3549
$ws.Range("J58").Value = 1749.75
$ws.Range("K58").Value = 1236.3549
$ws.Range("L58").Value = 1749.75
$ws.Range("M58").Value = -1033.3549
$ws.Range("N58").Value = -2155.75
$ws.Range("H86").Value = 2133
$ws.Range("J86").Value = 2199
$ws.Range("L86").Value = 2199
$ws.Range("N86").Value = -4445
$ws.Range("H89").Value = 2133
$ws.Range("J89").Value = 2199
$ws.Range("L89").Value = 10995
$ws.Range("N89").Value = -22227
$ws.Range("H134").Value = 1377.9487
$ws.Range("I134").Value = 1176.9656
$ws.Range("J134").Value = 1960.8
$ws.Range("K134").Value = 3530.8968
$ws.Range("L134").Value = 5882.4
$ws.Range("M134").Value = -995.8968
$ws.Range("N134").Value = -10952.4
$ws.Range("H136").Value = 1295.0286
$ws.Range("I136").Value = 1236.3549
$ws.Range("J136").Value = 1749.75
$ws.Range("K136").Value = 3709.0647
$ws.Range("L136").Value = 5249.25
$ws.Range("M136").Value = -1159.0647
$ws.Range("N136").Value = -10349.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 69.833336
$ws.Range("I14").Value = 69.833336
$ws.Range("K14").Value = 209.500008
$ws.Range("M14").Value = -36.50000800000001
$ws.Range("H80").Value = 5308.3335
$ws.Range("I80").Value = 16000
$ws.Range("J80").Value = 3170
$ws.Range("K80").Value = 48000
$ws.Range("L80").Value = 9510
$ws.Range("M80").Value = -47064
$ws.Range("N80").Value = -11382
$ws.Range("H83").Value = 5308.3335
$ws.Range("I83").Value = 16000
$ws.Range("J83").Value = 3170
$ws.Range("K83").Value = 144000
$ws.Range("L83").Value = 28530
$ws.Range("M83").Value = -139320
$ws.Range("N83").Value = -37890
$ws.Range("H107").Value = 461.55173
$ws.Range("I107").Value = 546.0833
$ws.Range("J107").Value = 401.88235
$ws.Range("K107").Value = 1638.2499
$ws.Range("L107").Value = 1205.64705
$ws.Range("M107").Value = 281.7501
$ws.Range("N107").Value = -5045.64705
$ws.Range("H122").Value = 709.2857
$ws.Range("I122").Value = 475.9091
$ws.Range("K122").Value = 4283.1819
$ws.Range("M122").Value = -1833.1819
$ws.Range("H131").Value = 880.74
$ws.Range("I131").Value = 575
$ws.Range("J131").Value = 886.9796
$ws.Range("K131").Value = 1725
$ws.Range("L131").Value = 2660.9388
$ws.Range("M131").Value = 3315
$ws.Range("N131").Value = -12740.9388

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6244.9443
$ws.Range("I70").Value = 5083.3335
$ws.Range("J70").Value = 6825.75
$ws.Range("K70").Value = 5083.3335
$ws.Range("L70").Value = 6825.75
$ws.Range("M70").Value = -4813.3335
$ws.Range("N70").Value = -7365.75
$ws.Range("H73").Value = 6244.9443
$ws.Range("I73").Value = 5083.3335
$ws.Range("J73").Value = 6825.75
$ws.Range("K73").Value = 5083.3335
$ws.Range("L73").Value = 6825.75
$ws.Range("M73").Value = -4147.3335
$ws.Range("N73").Value = -8697.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 71995.8
$ws.Range("J133").Value = 71995.8
$ws.Range("L133").Value = 71995.8
$ws.Range("N133").Value = -77055.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 12492.5
$ws.Range("I52").Value = 9990
$ws.Range("J52").Value = 20000
$ws.Range("K52").Value = 9990
$ws.Range("L52").Value = 20000
$ws.Range("M52").Value = -9764
$ws.Range("N52").Value = -20452
$ws.Range("H113").Value = 1122.2222
$ws.Range("I113").Value = 782.8333
$ws.Range("J113").Value = 1801
$ws.Range("K113").Value = 2348.4999
$ws.Range("L113").Value = 5403
$ws.Range("M113").Value = -178.4998999999998
$ws.Range("N113").Value = -9743
